# update after third manuscript draft
#
# - removes the two "Cochrane defines..." scratch notes (and their blank
#   separator row) from the top of the notes_overall sheet
# - replaces the old two-run (partly bold/red) CORIMUNO note with a plain
#   successor note now that the discrepancy has been resolved
# - moves the active tab / selection from the data sheet to notes_overall

$wb = $excel.ActiveWorkbook

$notes = $wb.Worksheets.Item("notes_overall")

# Drop the two Cochrane definition notes plus the blank row that used to
# separate them from the REMAP-CAP notes below - everything else shifts up.
$notes.Rows("1:3").Delete() | Out-Null

# The CORIMUNO discrepancy got resolved; swap in the follow-up note (plain
# text, no more bold red "(CHECK THIS)" run) in its place.
$notes.Range("A6").Value = "CORIMUNO death at D28 (Abstract/Table 2/eTable8 )"

# The data sheet is no longer the selected tab; update its lingering selection
# first (selecting it would otherwise re-activate it as the last-touched sheet).
$data = $wb.Worksheets.Item("data")
$data.Range("H5").Select() | Out-Null

# notes_overall becomes the active sheet/tab, with a fresh selection.
$notes.Activate()
$notes.Range("G8").Select() | Out-Null
